# Auto-generated: updates crypto price (D) and 1h volume % (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''61.985.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '''2.409.91'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.90%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''560.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.90%  '
$ws.Range('D6').Value = '''138.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.67%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '''0.586'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.99%  '
$ws.Range('D9').Value = '''2.407.91'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('E10').Value = '  +3.36%  '
$ws.Range('D11').Value = '''5.72'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.99%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  +3.98%  '
$ws.Range('D14').Value = '''25.69'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.55%  '
$ws.Range('D15').Value = '''2.842.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.94%  '
$ws.Range('D16').Value = '''61.918.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.78%  '
$ws.Range('E17').Value = '  +5.21%  '
$ws.Range('D18').Value = '''2.412.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.21%  '
$ws.Range('D19').Value = '''11.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.39%  '
$ws.Range('D20').Value = '''344.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.72%  '
$ws.Range('D21').Value = '''4.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('E22').Value = '  +3.30%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '''65.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').Value = '''0.172'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''1.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +14.28%  '
$ws.Range('E29').Value = '  +15.96%  '
$ws.Range('E30').Value = '  +4.17%  '
$ws.Range('D31').Value = '''0.0₃0781'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.03%  '
$ws.Range('E32').Value = '  +7.45%  '
$ws.Range('D33').Value = '''170.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.49%  '
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('E35').Value = '  +2.25%  '
$ws.Range('D36').Value = '''378.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +17.39%  '
$ws.Range('D37').Value = '''18.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.09%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Value = '''1.66'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.19%  '
$ws.Range('D42').Value = '''39.14'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.23%  '
$ws.Range('D43').Value = '''144.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.74%  '
$ws.Range('E44').Value = '  +4.88%  '
$ws.Range('D45').Value = '''20.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.06%  '
$ws.Range('E46').Value = '  +6.67%  '
$ws.Range('E47').Value = '  +2.09%  '
$ws.Range('D48').Value = '''0.586'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.63%  '
$ws.Range('D49').Value = '''17.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.34%  '
$ws.Range('E50').Value = '  +3.79%  '
$ws.Range('E51').Value = '  +1.29%  '
